$d = $word.ActiveDocument
$d.Content.Find.Execute("(11:00-", $true, $false, $false, $false, $false,
                         $true, 1, $false, "(11:00-13:30) (17:06-", 2)
